$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRA_data")

$oldText = "DY2012"
$newText = "DY20-12"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Text
    if ($valA.Contains($oldText)) {
        $cellA.Value = $valA.Replace($oldText, $newText)
    }

    $valB = $cellB.Text
    if ($valB.Contains($oldText)) {
        $cellB.Value = $valB.Replace($oldText, $newText)
    }
}
